$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "prueba" text to A2 and A3 (both reference the same new shared string)
$ws.Range("A2").Value = "prueba"
$ws.Range("A3").Value = "prueba"

# A4 is left blank but formatted with an underline font (no fill/border)
$ws.Range("A4").Font.Underline = $true

# Put the selection/active cell on A4, matching the saved view state
$ws.Range("A4").Select()
